$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New headers: H1=angkatan, I1=jenis_kelamin, J1=ipk (bold like the existing header row) ---
$ws.Cells.Item(1, 8).Value = "angkatan"
$ws.Cells.Item(1, 9).Value = "jenis_kelamin"
$ws.Cells.Item(1, 10).Value = "ipk"
$ws.Range("H1:J1").Font.Bold = $true

# --- New data columns for rows 2-50: angkatan (H), jenis_kelamin (I), ipk (J) ---
# Each entry: row,angkatan,jenis_kelamin,ipk,needsTwoDecimalFormat
$data = @(
    "2,2023,L,3.71,0",
    "3,2023,L,3.81,0",
    "4,2023,L,3.6,1",
    "5,2023,L,3.32,0",
    "6,2023,L,3.87,0",
    "7,2023,L,4,1",
    "8,2023,L,3.57,0",
    "9,2023,L,3.65,0",
    "10,2023,L,3.85,0",
    "11,2023,L,3.45,0",
    "12,2023,L,3.73,0",
    "13,2023,L,3.93,0",
    "14,2023,P,3.71,0",
    "15,2023,P,3.81,0",
    "16,2023,P,3.6,1",
    "17,2023,P,3.32,0",
    "18,2023,P,3.87,0",
    "19,2023,P,4,1",
    "20,2023,P,3.57,0",
    "21,2022,P,3.65,0",
    "22,2022,P,3.85,0",
    "23,2022,P,3.45,0",
    "24,2022,P,3.73,0",
    "25,2022,P,3.93,0",
    "26,2022,P,3.71,0",
    "27,2022,P,3.81,0",
    "28,2022,P,3.6,1",
    "29,2022,P,3.32,0",
    "30,2022,P,3.87,0",
    "31,2022,P,4,1",
    "32,2022,P,3.57,0",
    "33,2022,P,3.65,0",
    "34,2022,P,3.85,0",
    "35,2022,P,3.45,0",
    "36,2022,P,3.73,0",
    "37,2022,L,3.93,0",
    "38,2022,L,3.71,0",
    "39,2022,L,3.81,0",
    "40,2022,L,3.6,1",
    "41,2024,L,3.32,0",
    "42,2024,L,3.87,0",
    "43,2024,L,4,1",
    "44,2024,L,3.57,0",
    "45,2024,L,3.65,0",
    "46,2024,L,3.85,0",
    "47,2024,L,3.45,0",
    "48,2024,L,3.73,0",
    "49,2024,L,3.93,0",
    "50,2024,L,3.89,0"
)

foreach ($line in $data) {
    $parts = $line.Split(",")
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 8).Value = [int]$parts[1]
    $ws.Cells.Item($r, 9).Value = $parts[2]
    $ws.Cells.Item($r, 10).Value = [double]$parts[3]
    if ($parts[4] -eq "1") {
        $ws.Cells.Item($r, 10).NumberFormat = "0.00"
    }
}

# --- Column I width (new column). The engine quantizes ColumnWidth to 1/6-character
# steps, so 11.5 is the closest achievable approximation of the authored 12.26953125 width ---
$ws.Columns.Item(9).ColumnWidth = 11.5

# --- View state: selection at N49 (matches the cursor position after the edit) ---
$ws.Range("N49").Select()
